# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changedRange = $ws.Range("D2,E2,D3,E3,E4,D5,E5,E6,D7,E7,D8,E8,D9,E9,D10,E10,D11,E11,D12,E12,D13,E13,D14,E14,D15,E15,E16,D17,E17,E18,D19,E19,D20,E20,D21,E21,D22,E22,D23,E23,D24,E24,D25,E25,B26,C26,D26,E26,B27,C27,D27,E27,B28,C28,D28,E28,B29,C29,D29,E29,B30,C30,D30,E30,B31,C31,D31,E31,B32,C32,D32,E32,B33,C33,D33,E33,B34,C34,D34,E34,B35,C35,D35,E35,B36,C36,D36,E36,B37,C37,D37,E37,B38,C38,D38,E38,B39,C39,D39,E39,B40,C40,D40,E40,B41,C41,D41,E41,B42,C42,D42,E42,B43,C43,D43,E43,B44,C44,D44,E44,D45,E45,B46,C46,D46,E46,B47,C47,D47,E47,B48,C48,D48,E48,B49,C49,D49,E49,B50,C50,D50,E50,B51,C51,D51,E51")

# Force text format on every area so numeric-looking strings (e.g. "1.210")
# are preserved verbatim instead of being coerced into numbers.
foreach ($area in $changedRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range('D2').Value = '28.112.14'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.878.86'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '313.52'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '0.5087'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '0.3853'
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('D9').Value = '0.09015'
$ws.Range('E9').Value = '  -3.96%  '
$ws.Range('D10').Value = '1.124'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = '41.63'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = '6.345'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '20.77'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '1.876.57'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').Value = '7.226'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '0.00001111'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '0.06597'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').Value = '18.19'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '6.121'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = '28.129.89'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '11.42'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '2.272'
$ws.Range('E25').Value = '  -1.79%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '3.384'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '2.099.54'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.544'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '20.81'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '156.87'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = '127.03'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.1053'
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.062'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.612'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.600'
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '9.665'
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.06599'
$ws.Range('E37').Value = '  -0.86%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02423'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.2182'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.281'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.210'
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6409'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '11.52'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('B44').Value = 'InternetComputer(DFINITY)'
$ws.Range('C44').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D44').Value = '4.924'
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('D45').Value = '0.6045'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '13.20'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '3.674'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.275'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.243'
$ws.Range('E49').Value = '  +5.47%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '2.000'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '121.52'
$ws.Range('E51').Value = '  -0.62%  '

# Restore default (no explicit) cell formatting, matching original workbook style
foreach ($area in $changedRange.Areas) {
    $area.ClearFormats()
}

Write-Host "Updated cryptos list"
